# Apply a copy-editing pass over the resume: normalize "spaced hyphen"
# dashes to proper em/en dashes, and a handful of small wording tweaks.

$d = $word.ActiveDocument

function Replace-Text {
    param(
        [string]$find,
        [string]$replace
    )
    $null = $d.Content.Find.Execute(
        $find,        # FindText
        $true,        # MatchCase
        $false,       # MatchWholeWord
        $false,       # MatchWildcards
        $false,       # MatchSoundsLike
        $false,       # MatchAllWordForms
        $true,        # Forward
        1,            # Wrap (wdFindContinue)
        $false,       # Format
        $replace,     # ReplaceWith
        2             # Replace (wdReplaceAll)
    )
}

# Master's thesis sentence dash -> em dash
Replace-Text "I used Python for my Master’s thesis -" "I used Python for my Master’s thesis —"

# "changelogs" -> "change logs"; "security impact analyses" -> "security impact reports"
Replace-Text "by monitoring releases, reviewing changelogs, integrating version updates, submitting security impact analyses, and training co-workers on modern trends and technologies." "by monitoring releases, reviewing change logs, integrating version updates, submitting security impact reports, and training co-workers on modern trends and technologies."

# "coding ability" sentence dash -> em dash
Replace-Text "for examples of my coding ability - like my PowerShell module," "for examples of my coding ability — like my PowerShell module,"

# "technical evolutions" -> "technical procedures"
Replace-Text "Led teams through complex technical evolutions, during maintenance and casualty situations, maximizing efficiency while maintaining a safe environment." "Led teams through complex technical procedures, during maintenance and casualty situations, maximizing efficiency while maintaining a safe environment."

# Date ranges: hyphen -> en dash
Replace-Text "03/20 - Present," "03/20 – Present,"
Replace-Text "11/15 - 03/20," "11/15 – 03/20,"
Replace-Text "06/18 - 07/19," "06/18 – 07/19,"
Replace-Text "08/15 - 11/15," "08/15 – 11/15,"
Replace-Text "07/14 - 08/15," "07/14 – 08/15,"
Replace-Text "07/10 - 07/14," "07/10 – 07/14,"

# "decision making" -> "decision-making"
Replace-Text "enabled collaboration and senior leader strategic decision making." "enabled collaboration and senior leader strategic decision-making."

# "Web Apps Must be Sustainable" sentence dash -> em dash
Replace-Text "Web Apps Must be Sustainable - How to leverage common patterns to enable flexibility" "Web Apps Must be Sustainable — How to leverage common patterns to enable flexibility"

# 2019 line: hyphen -> em dash
Replace-Text "2019 - Completed OSCP Penetration with Kali course. I did not pass the 24hr OSCP exam, but I learned a lot in the process." "2019 — Completed OSCP Penetration with Kali course. I did not pass the 24hr OSCP exam, but I learned a lot in the process."

# 2014 line: hyphen -> em dash
Replace-Text "2014 - Present Navy Selected Ready Reserve member (Commander, O-5)" "2014 — Present Navy Selected Ready Reserve member (Commander, O-5)"

# 2012/2011 date ranges: hyphen -> en dash
Replace-Text "2012 - 2013 Tutored online with tutor.com in Calculus and Algebra" "2012 – 2013 Tutored online with tutor.com in Calculus and Algebra"
Replace-Text "2011 - 2013 Tutored local high school students in Geometry and Trigonometry" "2011 – 2013 Tutored local high school students in Geometry and Trigonometry"

Write-Output "edits applied"
